# "update file with jgit" — the rules sheet had its "From" greeting label
# for rule R10 (cell E8) changed from "Good Morning" to "GIT UPDATE".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

# Leave the selection on the edited cell, matching the saved sheetView state.
$ws.Range("E8").Select() | Out-Null
